# Automatic update of files.
# Update the "Förändrad" (column C) date for all existing data rows (2-27)
# from 45202 (2023-10-03) to 45203 (2023-10-04), and append two new data
# rows (28, 29) for newly reported cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump column C (Förändrad) from 45202 to 45203 for rows 2..27 ---
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}

# --- 2. Append new row 28: A 47008-2023 ---
$row = 28
$ws.Cells.Item($row, 1).Value = "A 47008-2023"
$ws.Cells.Item($row, 2).Value = 45201
$ws.Cells.Item($row, 3).Value = 45203
$ws.Cells.Item($row, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($row, 5).Value = "HAMMARÖ"
$ws.Cells.Item($row, 6).Value = "Kommuner"
$ws.Cells.Item($row, 7).Value = 3.3
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = "'"

# --- 3. Append new row 29: A 46982-2023 ---
$row = 29
$ws.Cells.Item($row, 1).Value = "A 46982-2023"
$ws.Cells.Item($row, 2).Value = 45201
$ws.Cells.Item($row, 3).Value = 45203
$ws.Cells.Item($row, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($row, 5).Value = "HAMMARÖ"
$ws.Cells.Item($row, 6).Value = "Kommuner"
$ws.Cells.Item($row, 7).Value = 3.8
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = "'"

# --- 4. Apply the same number format (date) used by B/C columns to the new rows ---
$ws.Range("B28:C29").NumberFormat = $ws.Range("B27:C27").NumberFormat

# --- 5. Apply the same style (wrap text, empty inline string) used by column R
#        to the new rows' R cells. A plain Value/Style assignment on an
#        empty string collapses back to a blank cell in this engine, so we
#        go through Copy / PasteSpecial(xlPasteFormats) from the existing
#        R27 cell (which already carries the right style) to stamp the
#        identical style index onto R28/R29 without disturbing their content.
$ws.Range("R27").Copy()
$ws.Range("R28").PasteSpecial(-4122)
$ws.Range("R27").Copy()
$ws.Range("R29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Match row heights: row 27/28 get explicit 15pt custom height, row 29 keeps default ---
$ws.Rows.Item(27).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 15
